$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1349.9445
$ws.Range("I41").Value = 1416.1111
$ws.Range("J41").Value = 1283.7778
$ws.Range("K41").Value = 1416.1111
$ws.Range("L41").Value = 1283.7778
$ws.Range("M41").Value = -976.1111000000001
$ws.Range("N41").Value = -2163.7778

$ws.Range("H103").Value = 1909.8
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1909.8
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 5729.4
$ws.Range("N103").Value = -6901.4
$ws.Range("M103").ClearContents()

$ws.Range("H138").Value = 1374.52
$ws.Range("I138").Value = 874.4737
$ws.Range("J138").Value = 1681
$ws.Range("K138").Value = 2623.4211
$ws.Range("L138").Value = 5043
$ws.Range("M138").Value = 2516.5789
$ws.Range("N138").Value = -15323

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3660
$ws.Range("I32").Value = 3289.6956
$ws.Range("J32").Value = 9338
$ws.Range("K32").Value = 3289.6956
$ws.Range("L32").Value = 9338
$ws.Range("M32").Value = -3002.6956
$ws.Range("N32").Value = -9912

$ws.Range("H74").Value = 1555.7142
$ws.Range("I74").Value = 974.8570999999999
$ws.Range("J74").Value = 2136.5715
$ws.Range("K74").Value = 974.8570999999999
$ws.Range("L74").Value = 2136.5715
$ws.Range("M74").Value = -100.8570999999999
$ws.Range("N74").Value = -3884.5715

$ws.Range("H77").Value = 1555.7142
$ws.Range("I77").Value = 974.8570999999999
$ws.Range("J77").Value = 2136.5715
$ws.Range("K77").Value = 4874.2855
$ws.Range("L77").Value = 10682.8575
$ws.Range("M77").Value = -506.2855
$ws.Range("N77").Value = -19418.8575

$ws.Range("H110").Value = 1238.8422
$ws.Range("I110").Value = 727.2727
$ws.Range("J110").Value = 1942.25
$ws.Range("K110").Value = 727.2727
$ws.Range("L110").Value = 1942.25
$ws.Range("M110").Value = 1317.7273
$ws.Range("N110").Value = -6032.25

$ws.Range("H122").Value = 1572.1
$ws.Range("I122").Value = 1634.1111
$ws.Range("J122").Value = 1014
$ws.Range("K122").Value = 4902.3333
$ws.Range("L122").Value = 3042
$ws.Range("M122").Value = -2452.3333
$ws.Range("N122").Value = -7942

$ws.Range("H133").Value = 28672.723
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 28672.723
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 28672.723
$ws.Range("N133").Value = -33732.723

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 15625839
$ws.Range("I94").Value = 19231526
$ws.Range("J94").Value = 1196.3334
$ws.Range("K94").Value = 19231526
$ws.Range("L94").Value = 1196.3334
$ws.Range("M94").Value = -19231075

$ws.Range("H99").Value = 50001570
$ws.Range("I99").Value = 55557076
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 55557076
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -55555578

$ws.Range("H105").Value = 168316510
$ws.Range("I105").Value = 201979500
$ws.Range("J105").Value = 1555
$ws.Range("K105").Value = 201979500
$ws.Range("L105").Value = 1555
$ws.Range("M105").Value = -201977753
$ws.Range("N105").Value = -5049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -387

$ws.Range("H16").Value = 37037890
$ws.Range("I16").Value = 43479090
$ws.Range("J16").Value = 970
$ws.Range("K16").Value = 43479090
$ws.Range("L16").Value = 970
$ws.Range("M16").Value = -43478803
$ws.Range("N16").Value = -1544

$ws.Range("H41").Value = 22666.666
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 22666.666
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 22666.666
$ws.Range("N41").Value = -23522.666
$ws.Range("M41").ClearContents()

$ws.Range("H58").Value = 1356.2941
$ws.Range("I58").Value = 1162.875
$ws.Range("J58").Value = 1820.5
$ws.Range("K58").Value = 1162.875
$ws.Range("L58").Value = 1820.5
$ws.Range("M58").Value = -959.875

$ws.Range("H112").Value = 37714.285
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 37714.285
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 37714.285
$ws.Range("N112").Value = -40668.285

$ws.Range("H113").Value = 37037890
$ws.Range("I113").Value = 43479090
$ws.Range("J113").Value = 970
$ws.Range("K113").Value = 43479090
$ws.Range("L113").Value = 970
$ws.Range("M113").Value = -43476920
$ws.Range("N113").Value = -5310

$ws.Range("H132").Value = 7959.65
$ws.Range("I132").Value = 16733.143
$ws.Range("J132").Value = 3235.4614
$ws.Range("K132").Value = 50199.429
$ws.Range("L132").Value = 9706.3842
$ws.Range("M132").Value = -47669.429
$ws.Range("N132").Value = -14766.3842

$ws.Range("H134").Value = 1709.9697
$ws.Range("I134").Value = 1686.1923
$ws.Range("J134").Value = 1798.2858
$ws.Range("K134").Value = 5058.5769
$ws.Range("L134").Value = 5394.857400000001
$ws.Range("M134").Value = -2523.5769
$ws.Range("N134").Value = -10464.8574

$ws.Range("H136").Value = 1356.2941
$ws.Range("I136").Value = 1162.875
$ws.Range("J136").Value = 1820.5
$ws.Range("K136").Value = 3488.625
$ws.Range("L136").Value = 5461.5
$ws.Range("M136").Value = -938.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 81.20689400000001
$ws.Range("I12").Value = 124
$ws.Range("J12").Value = 61.95
$ws.Range("K12").Value = 372
$ws.Range("L12").Value = 185.85
$ws.Range("M12").Value = -199
$ws.Range("N12").Value = -531.85

$ws.Range("H98").Value = 200
$ws.Range("I98").Value = 150
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 450
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = 1048
$ws.Range("N98").Value = -3896

$ws.Range("H131").Value = 16394680
$ws.Range("I131").Value = 250000350
$ws.Range("J131").Value = 1300.0702
$ws.Range("K131").Value = 750001050
$ws.Range("L131").Value = 3900.2106
$ws.Range("M131").Value = -749996010
$ws.Range("N131").Value = -13980.2106

$ws.Range("H132").Value = 1121.6154
$ws.Range("I132").Value = 682
$ws.Range("J132").Value = 1317
$ws.Range("K132").Value = 6138
$ws.Range("L132").Value = 11853
$ws.Range("M132").Value = -3608
$ws.Range("N132").Value = -16913

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1950
$ws.Range("I4").Value = 900
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -788
$ws.Range("N4").Value = -3224

$ws.Range("H97").Value = 612.1429000000001
$ws.Range("I97").Value = 630
$ws.Range("J97").Value = 567.5
$ws.Range("K97").Value = 630
$ws.Range("L97").Value = 567.5
$ws.Range("M97").Value = -134
$ws.Range("N97").Value = -1559.5

$ws.Range("H122").Value = 1691
$ws.Range("I122").Value = 2044.1111
$ws.Range("J122").Value = 1055.4
$ws.Range("K122").Value = 6132.3333
$ws.Range("L122").Value = 3166.2
$ws.Range("M122").Value = -3682.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 547.6667
$ws.Range("I16").Value = 537.8
$ws.Range("J16").Value = 597
$ws.Range("K16").Value = 537.8
$ws.Range("L16").Value = 597
$ws.Range("M16").Value = -367.8
$ws.Range("N16").Value = -937

$ws.Range("H21").Value = 3500
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3500
$ws.Range("N21").Value = -3848

$ws.Range("H40").Value = 3012.5
$ws.Range("I40").Value = 2867.6667
$ws.Range("J40").Value = 3099.4
$ws.Range("K40").Value = 2867.6667
$ws.Range("L40").Value = 3099.4
$ws.Range("M40").Value = -2731.6667
$ws.Range("N40").Value = -3371.4

$ws.Range("H122").Value = 20843234
$ws.Range("I122").Value = 20843234
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 62529702
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -62527252
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 976.1667
$ws.Range("I96").Value = 794
$ws.Range("J96").Value = 2980
$ws.Range("K96").Value = 794
$ws.Range("L96").Value = 2980
$ws.Range("M96").Value = 579
$ws.Range("N96").Value = -5726

$ws.Range("H107").Value = 467.875
$ws.Range("I107").Value = 448
$ws.Range("J107").Value = 501
$ws.Range("K107").Value = 1344
$ws.Range("L107").Value = 1503
$ws.Range("M107").Value = 576
$ws.Range("N107").Value = -5343

$ws.Range("H112").Value = 28933.334
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 28933.334
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 28933.334
$ws.Range("N112").Value = -31887.334

$ws.Range("H113").Value = 414.0909
$ws.Range("I113").Value = 230.28572
$ws.Range("J113").Value = 735.75
$ws.Range("K113").Value = 690.85716
$ws.Range("L113").Value = 2207.25
$ws.Range("M113").Value = 1479.14284
$ws.Range("N113").Value = -6547.25

$ws.Range("H126").Value = 66667880
$ws.Range("I126").Value = 100001050
$ws.Range("J126").Value = 1542
$ws.Range("K126").Value = 300003150
$ws.Range("L126").Value = 4626
$ws.Range("M126").Value = -300000680
$ws.Range("N126").Value = -9566
